$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Reserva")
$c = $ws.Range("A47")
$c.Formula = "=TEXT(TRUE,""General"")"
Write-Host "val:" $c.Value
Write-Host "type via GetCellType equivalent"
$c.Copy()
$c.PasteSpecial(-4163)  # xlPasteValues
